$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.75
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 2.9
$ws.Range("Z5").Value = 11
$ws.Range("AA5").Value = 11
$ws.Range("G11").Value = 1.6
$ws.Range("H11").Value = 3.9
$ws.Range("I11").Value = 5.75
$ws.Range("J11").Value = 2.25
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("U11").Value = 1.5
$ws.Range("V11").Value = 2.5
$ws.Range("AB11").Value = 11
$ws.Range("AE11").Value = 8
$ws.Range("AF11").Value = 7.5
$ws.Range("G12").Value = 2.32
$ws.Range("H12").Value = 2.75
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 1.88
$ws.Range("L12").Value = 4.1
$ws.Range("M12").Value = 1.13
$ws.Range("N12").Value = 5.2
$ws.Range("O12").Value = 1.52
$ws.Range("P12").Value = 2.37
$ws.Range("Q12").Value = 2.52
$ws.Range("R12").Value = 1.47
$ws.Range("S12").Value = 4.5
$ws.Range("T12").Value = 1.17
$ws.Range("U12").Value = 1.55
$ws.Range("V12").Value = 2.3
$ws.Range("W12").Value = 2.05
$ws.Range("X12").Value = 1.7
$ws.Range("Y12").Value = 5.8
$ws.Range("AA12").Value = 9.5
$ws.Range("AE12").Value = 5.2
$ws.Range("AG12").Value = 16.5
$ws.Range("AJ12").Value = 7.8
$ws.Range("AK12").Value = 17
$ws.Range("AL12").Value = 12.5
$ws.Range("AN12").Value = 40
$ws.Range("AO12").Value = 55
$ws.Range("G13").Value = 2.67
$ws.Range("H13").Value = 2.7
$ws.Range("I13").Value = 3
$ws.Range("L13").Value = 3.6
$ws.Range("O13").Value = 1.55
$ws.Range("P13").Value = 2.3
$ws.Range("Q13").Value = 2.6
$ws.Range("R13").Value = 1.44
$ws.Range("S13").Value = 4.65
$ws.Range("Y13").Value = 6.2
$ws.Range("AC13").Value = 29
$ws.Range("AD13").Value = 50
$ws.Range("AG13").Value = 17
$ws.Range("AO13").Value = 45
$ws.Range("G14").Value = 2.45
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 3.1
$ws.Range("O14").Value = 1.5
$ws.Range("P14").Value = 2.63
$ws.Range("U14").Value = 1.57
$ws.Range("V14").Value = 2.25
$ws.Range("W14").Value = 2.1
$ws.Range("X14").Value = 1.67
$ws.Range("Y14").Value = 6.5
$ws.Range("AB14").Value = 23
$ws.Range("AE14").Value = 6.5
$ws.Range("AI14").Value = 501
$ws.Range("AL14").Value = 12
$ws.Range("AM14").Value = 34
$ws.Range("AN14").Value = 29
$ws.Range("AP14").Value = 1.85
$ws.Range("AQ14").Value = 2
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 13
$ws.Range("I17").Value = 6
$ws.Range("L17").Value = 6
$ws.Range("Q17").Value = 1.88
$ws.Range("R17").Value = 1.98
$ws.Range("AB17").Value = 11
$ws.Range("AD17").Value = 26
$ws.Range("AF17").Value = 7.5
$ws.Range("G18").Value = 3.1
$ws.Range("I18").Value = 2.25
$ws.Range("AF18").Value = 6
$ws.Range("AN18").Value = 21
$ws.Range("G19").Value = 2.4
$ws.Range("H19").Value = 3.25
$ws.Range("I19").Value = 2.88
$ws.Range("J19").Value = 3.2
$ws.Range("K19").Value = 2.05
$ws.Range("N19").Value = 10
$ws.Range("Z19").Value = 12
$ws.Range("AA19").Value = 10
$ws.Range("AE19").Value = 9
$ws.Range("AF19").Value = 6
$ws.Range("AJ19").Value = 8.5
$ws.Range("AK19").Value = 13
$ws.Range("G21").Value = 2.4
$ws.Range("M21").Value = 1.05
$ws.Range("N21").Value = 11
$ws.Range("Q21").Value = 1.95
$ws.Range("R21").Value = 1.9
$ws.Range("S21").Value = 3.25
$ws.Range("T21").Value = 1.33
$ws.Range("W21").Value = 1.7
$ws.Range("X21").Value = 2.05
$ws.Range("AA21").Value = 9.5
$ws.Range("AC21").Value = 19
$ws.Range("AE21").Value = 11
$ws.Range("AJ21").Value = 9.5
$ws.Range("G23").Value = 1.44
$ws.Range("H23").Value = 4.33
$ws.Range("Q23").Value = 2.03
$ws.Range("R23").Value = 1.83
$ws.Range("W23").Value = 2.1
$ws.Range("X23").Value = 1.67
$ws.Range("Z23").Value = 6.5
$ws.Range("AE23").Value = 10
$ws.Range("AG23").Value = 21
$ws.Range("G24").Value = 1.9
$ws.Range("H24").Value = 3.25
$ws.Range("K24").Value = 1.95
$ws.Range("M24").Value = 1.1
$ws.Range("N24").Value = 7
$ws.Range("Q24").Value = 2.5
$ws.Range("R24").Value = 1.5
$ws.Range("S24").Value = 5
$ws.Range("T24").Value = 1.17
$ws.Range("AA24").Value = 9.5
$ws.Range("AD24").Value = 41
$ws.Range("AI24").Value = 501
$ws.Range("AJ24").Value = 9.5
$ws.Range("AP24").Value = 1.88
$ws.Range("AQ24").Value = 1.98
$ws.Range("Q25").Value = 1.93
$ws.Range("R25").Value = 1.93
$ws.Range("S25").Value = 3.25
$ws.Range("T25").Value = 1.33
$ws.Range("M26").Value = 1.07
$ws.Range("N26").Value = 9
$ws.Range("O26").Value = 1.36
$ws.Range("P26").Value = 3
$ws.Range("Q26").Value = 2.2
$ws.Range("R26").Value = 1.65
$ws.Range("S26").Value = 4
$ws.Range("T26").Value = 1.22
$ws.Range("J28").Value = 2.75
$ws.Range("K28").Value = 2.12
$ws.Range("O28").Value = 1.2
$ws.Range("P28").Value = 3.6
$ws.Range("Q28").Value = 1.62
$ws.Range("R28").Value = 2.05
$ws.Range("S28").Value = 2.42
$ws.Range("T28").Value = 1.44
$ws.Range("W28").Value = 1.5
$ws.Range("X28").Value = 2.27
$ws.Range("Y28").Value = 10
$ws.Range("Z28").Value = 13
$ws.Range("AB28").Value = 24
$ws.Range("AC28").Value = 16
$ws.Range("AD28").Value = 21
$ws.Range("AE28").Value = 12.5
$ws.Range("AF28").Value = 6.6
$ws.Range("AG28").Value = 11
$ws.Range("AH28").Value = 40
$ws.Range("AI28").Value = 250
$ws.Range("AJ28").Value = 11.75
$ws.Range("AO28").Value = 25
$ws.Range("G30").Value = 2.45
$ws.Range("I30").Value = 3
$ws.Range("J30").Value = 3.1
$ws.Range("K30").Value = 2.05
$ws.Range("L30").Value = 3.6
$ws.Range("N30").Value = 8.5
$ws.Range("Q30").Value = 2.1
$ws.Range("R30").Value = 1.7
$ws.Range("Y30").Value = 7.5
$ws.Range("Z30").Value = 11
$ws.Range("AA30").Value = 9.5
$ws.Range("AE30").Value = 8.5
$ws.Range("AI30").Value = 301
$ws.Range("AK30").Value = 15
$ws.Range("AN30").Value = 26
$ws.Range("M31").Value = 1.08
$ws.Range("N31").Value = 8
$ws.Range("O31").Value = 1.4
$ws.Range("P31").Value = 2.75
$ws.Range("Q31").Value = 2.35
$ws.Range("R31").Value = 1.57
$ws.Range("S31").Value = 4.33
$ws.Range("T31").Value = 1.2
$ws.Range("AI31").Value = 1000
$ws.Range("AP31").Value = 1.75
$ws.Range("AQ31").Value = 2.05
$ws.Range("Q33").Value = 2.05
$ws.Range("R33").Value = 1.8
$ws.Range("S33").Value = 3.5
$ws.Range("T33").Value = 1.29
